$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Row 1 (0-based row 0)
$t.Cell(1,1).Range.Text = "814÷6="
$t.Cell(1,2).Range.Text = "411÷9="
$t.Cell(1,3).Range.Text = "470÷6="
$t.Cell(1,4).Range.Text = "542÷6="
$t.Cell(1,5).Range.Text = "302÷5="

# Row 5 (0-based row 4)
$t.Cell(5,1).Range.Text = "576÷6="
$t.Cell(5,2).Range.Text = "130÷3="
$t.Cell(5,3).Range.Text = "459÷7="
$t.Cell(5,4).Range.Text = "125÷4="
$t.Cell(5,5).Range.Text = "782÷8="

# Row 9 (0-based row 8)
$t.Cell(9,1).Range.Text = "493÷9="
$t.Cell(9,2).Range.Text = "895÷5="
$t.Cell(9,3).Range.Text = "162÷2="
$t.Cell(9,4).Range.Text = "703÷6="
$t.Cell(9,5).Range.Text = "833÷7="

# Row 13 (0-based row 12)
$t.Cell(13,1).Range.Text = "541÷8="
$t.Cell(13,2).Range.Text = "627÷3="
$t.Cell(13,3).Range.Text = "992÷2="
$t.Cell(13,4).Range.Text = "845÷3="
$t.Cell(13,5).Range.Text = "562÷6="

# Row 17 (0-based row 16)
$t.Cell(17,1).Range.Text = "859÷5="
$t.Cell(17,2).Range.Text = "710÷4="
$t.Cell(17,3).Range.Text = "320÷6="
$t.Cell(17,4).Range.Text = "644÷2="
$t.Cell(17,5).Range.Text = "940÷2="
